$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price / 1h volume-change figures, and fix a swapped
# WEMIXToken/Monero row pair (rows 41-42), per the Feb 9 2024 data refresh.

$ws.Range("D2").Value = '47.270.54'
$ws.Range("E2").Value = '  +4.37%  '
$ws.Range("D3").Value = '2.491.95'
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.01'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("E6").Value = '  +5.22%  '
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.25'
$ws.Range("E10").Value = '  +7.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.42'
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.20'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("D15").Value = '2.880.52'
$ws.Range("E15").Value = '  +2.70%  '
$ws.Range("D16").Value = '2.495.74'
$ws.Range("E16").Value = '  +4.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.856'
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").Value = '47.197.28'
$ws.Range("E18").Value = '  +4.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.96'
$ws.Range("E19").Value = '  +6.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.68'
$ws.Range("E20").Value = '  +5.55%  '
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.64'
$ws.Range("E22").Value = '  +2.41%  '
$ws.Range("E23").Value = '  +6.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '250.51'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.60'
$ws.Range("E25").Value = '  +4.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.19'
$ws.Range("E26").Value = '  +1.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  +4.27%  '
$ws.Range("E29").Value = '  -3.30%  '
$ws.Range("E30").Value = '  +11.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.08'
$ws.Range("E31").Value = '  +6.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.35'
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.50'
$ws.Range("E33").Value = '  +5.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.62'
$ws.Range("E34").Value = '  -2.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0791'
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.99'
$ws.Range("E37").Value = '  +6.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.68'
$ws.Range("E38").Value = '  +5.34%  '
$ws.Range("E39").Value = '  +4.01%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '121.97'
$ws.Range("E41").Value = '  -3.97%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.24'
$ws.Range("E42").Value = '  +1.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.27'
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("E44").Value = '  +3.05%  '
$ws.Range("D45").Value = '1.966.71'
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.02'
$ws.Range("E46").Value = '  +2.26%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.80'
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.29'
$ws.Range("E50").Value = '  +9.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.48'
$ws.Range("E51").Value = '  +4.01%  '
